# chore: adapt column header formatting to respective input file names
#
# Renames the header row from the generic "_old"/"_new" suffixes to the
# format-version-specific "_FV2410"/"_FV2504" suffixes, wraps the data range
# in a native Excel Table (ListObject) with autofilter, and freezes the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row (row 1) cells: "<name>_old" -> "<name>_FV2410",
#    "<name>_new" -> "<name>_FV2504". Column K ("diff") and everything else
#    is left untouched.
$newHeaders = @{
    1  = "Segmentname_FV2410"
    2  = "Segmentgruppe_FV2410"
    3  = "Segment_FV2410"
    4  = "Datenelement_FV2410"
    5  = "Segment ID_FV2410"
    6  = "Code_FV2410"
    7  = "Qualifier_FV2410"
    8  = "Beschreibung_FV2410"
    9  = "Bedingungsausdruck_FV2410"
    10 = "Bedingung_FV2410"
    12 = "Segmentname_FV2504"
    13 = "Segmentgruppe_FV2504"
    14 = "Segment_FV2504"
    15 = "Datenelement_FV2504"
    16 = "Segment ID_FV2504"
    17 = "Code_FV2504"
    18 = "Qualifier_FV2504"
    19 = "Beschreibung_FV2504"
    20 = "Bedingungsausdruck_FV2504"
    21 = "Bedingung_FV2504"
}

foreach ($col in $newHeaders.Keys) {
    $ws.Cells.Item(1, $col).Value2 = $newHeaders[$col]
}

# 2. Turn the used range into a proper Excel Table ("Table1") with an
#    autofilter, matching the table definition added at xl/tables/table1.xml.
$dataRange = $ws.Range("A1:U62")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# 3. Freeze the header row (split below row 1) so column headers stay
#    visible while scrolling through the 62 data rows.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
